$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old sheet had a two-row header (row1: column group labels in E/G/I/J/K,
# row2: unit labels in F/G/H/I/J/K) followed by 11 data rows (rows 3-13).
# The new layout collapses this into a single header row (row1) with a
# proper column per field, followed by the same 11 data rows shifted up by
# one (rows 2-12). Remove the old "units" row (row 2); this shifts the data
# rows up automatically and keeps their values/strings intact.
$ws.Rows.Item(2).Delete()

# Rebuild row 1 as the new header row.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# The A1:E1 labels are plain, unstyled cells - E1 inherited styling from the
# old layout (it used to hold the "(m3/s)" unit label), so reset it back to
# the default look.
$ws.Range("A1:E1").Style = "Normal"

# Match the font used elsewhere in the sheet (Arial 9) for the header cells
# that carry the new "(m3/s)/(MW1)/(MW2)/(GWh) ..." unit-style formatting.
$ws.Range("F1:K1").Font.Name = "Arial"
$ws.Range("F1:K1").Font.Size = 9

# Update the active selection to match the new focus (first data row).
$ws.Range("A2:K2").Select()
